$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2835.9033
$ws.Range("I132").Value = 2763.7666
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 8291.2998
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -5761.299800000001
$ws.Range("N132").Value = -20060

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 195.5
$ws.Range("I4").Value = 219.6
$ws.Range("K4").Value = 219.6
$ws.Range("M4").Value = -103.6

$ws.Range("H61").Value = 50001760
$ws.Range("I61").Value = 50001760
$ws.Range("K61").Value = 50001760
$ws.Range("M61").Value = -50001548

$ws.Range("H122").Value = 13336962
$ws.Range("I122").Value = 3410.6667
$ws.Range("K122").Value = 10232.0001
$ws.Range("M122").Value = -7782.000100000001

$ws.Range("H132").Value = 20865198
$ws.Range("I132").Value = 3043.1177
$ws.Range("J132").Value = 71530430
$ws.Range("K132").Value = 9129.3531
$ws.Range("L132").Value = 214591290
$ws.Range("M132").Value = -6599.3531
$ws.Range("N132").Value = -214596350

$ws.Range("H136").Value = 50001760
$ws.Range("I136").Value = 50001760
$ws.Range("K136").Value = 150005280
$ws.Range("M136").Value = -150002730

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 3728.8
$ws.Range("I99").Value = 2868.3333
$ws.Range("K99").Value = 2868.3333
$ws.Range("M99").Value = -1370.3333

$ws.Range("H105").Value = 6955.269
$ws.Range("I105").Value = 8509.125
$ws.Range("J105").Value = 4469.1
$ws.Range("K105").Value = 8509.125
$ws.Range("L105").Value = 4469.1
$ws.Range("M105").Value = -6762.125
$ws.Range("N105").Value = -7963.1

$ws.Range("H107").Value = 2767.4211
$ws.Range("I107").Value = 1398.8572
$ws.Range("J107").Value = 6599.4
$ws.Range("K107").Value = 1398.8572
$ws.Range("L107").Value = 6599.4
$ws.Range("M107").Value = 521.1428000000001
$ws.Range("N107").Value = -10439.4

$ws.Range("H109").Value = 120000
$ws.Range("J109").Value = 120000
$ws.Range("L109").Value = 120000
$ws.Range("N109").Value = -122774

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 10884.577
$ws.Range("I99").Value = 7374.75
$ws.Range("J99").Value = 11522.728
$ws.Range("K99").Value = 7374.75
$ws.Range("L99").Value = 11522.728
$ws.Range("M99").Value = -5876.75
$ws.Range("N99").Value = -14518.728

$ws.Range("H126").Value = 10884.577
$ws.Range("I126").Value = 7374.75
$ws.Range("J126").Value = 11522.728
$ws.Range("K126").Value = 22124.25
$ws.Range("L126").Value = 34568.18399999999
$ws.Range("M126").Value = -19654.25
$ws.Range("N126").Value = -39508.18399999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 12562.5
$ws.Range("I2").Value = 12562.5
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 75375
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -75262
$ws.Range("N2").ClearContents()

$ws.Range("H11").Value = 156.66667

$ws.Range("H17").Value = 340
$ws.Range("J17").Value = 360
$ws.Range("L17").Value = 1080
$ws.Range("N17").Value = -1418

$ws.Range("H18").Value = 2521.6667
$ws.Range("I18").Value = 1026
$ws.Range("K18").Value = 3078
$ws.Range("M18").Value = -2909

$ws.Range("H20").Value = 200
$ws.Range("I20").Value = 200
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 600
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = -373
$ws.Range("N20").ClearContents()

$ws.Range("H102").Value = 2750
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 8250
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -5816
$ws.Range("N102").ClearContents()

$ws.Range("H107").Value = 1942.7142
$ws.Range("I107").Value = 200
$ws.Range("J107").Value = 2233.1667
$ws.Range("K107").Value = 600
$ws.Range("L107").Value = 6699.500100000001
$ws.Range("M107").Value = 1320
$ws.Range("N107").Value = -10539.5001

$ws.Range("H112").Value = 10228.111
$ws.Range("I112").Value = 2513.25
$ws.Range("K112").Value = 7539.75
$ws.Range("M112").Value = -6431.75

$ws.Range("H123").Value = 2029
$ws.Range("J123").Value = 2029
$ws.Range("L123").Value = 6087
$ws.Range("N123").Value = -10987

$ws.Range("H131").Value = 5245.0625
$ws.Range("I131").Value = 1309.6
$ws.Range("J131").Value = 7033.909
$ws.Range("K131").Value = 3928.8
$ws.Range("L131").Value = 21101.727
$ws.Range("M131").Value = 1111.2
$ws.Range("N131").Value = -31181.727

$ws.Range("H133").Value = 20000
$ws.Range("J133").Value = 20000
$ws.Range("L133").Value = 60000
$ws.Range("N133").Value = -70120

$ws.Range("H136").Value = 2425.8
$ws.Range("I136").Value = 2425.8
$ws.Range("K136").Value = 7277.400000000001
$ws.Range("M136").Value = -2177.400000000001

$ws.Range("H137").Value = 2153.111
$ws.Range("J137").Value = 3029.8
$ws.Range("L137").Value = 9089.400000000001
$ws.Range("N137").Value = -19289.4

$ws.Range("H138").Value = 3000
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 3000
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 9000
$ws.Range("N138").Value = -19280
$ws.Range("M138").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 71361.664
$ws.Range("J62").Value = 70085
$ws.Range("L62").Value = 70085
$ws.Range("N62").Value = -71457

$ws.Range("H65").Value = 71361.664
$ws.Range("J65").Value = 70085
$ws.Range("L65").Value = 210255
$ws.Range("N65").Value = -217119

$ws.Range("H124").Value = 90333
$ws.Range("J124").Value = 90333
$ws.Range("L124").Value = 90333
$ws.Range("N124").Value = -100153

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2262.2222
$ws.Range("J22").Value = 3018.75
$ws.Range("L22").Value = 3018.75
$ws.Range("N22").Value = -3608.75

$ws.Range("H27").Value = 2262.2222
$ws.Range("J27").Value = 3018.75
$ws.Range("L27").Value = 3018.75
$ws.Range("N27").Value = -3232.75

$ws.Range("H46").Value = 1882
$ws.Range("I46").Value = 1019.5217
$ws.Range("J46").Value = 4086.111
$ws.Range("K46").Value = 1019.5217
$ws.Range("L46").Value = 4086.111
$ws.Range("M46").Value = -831.5217
$ws.Range("N46").Value = -4462.111

$ws.Range("H132").Value = 100003624
$ws.Range("J132").Value = 333339300
$ws.Range("L132").Value = 1000017900
$ws.Range("N132").Value = -1000022960

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H129").Value = 69680.164
$ws.Range("J129").Value = 69680.164
$ws.Range("L129").Value = 69680.164
$ws.Range("N129").Value = -79680.164
